$wb = $excel.ActiveWorkbook

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 18757512
$ws.Range("I74").Value = 37506524
$ws.Range("J74").Value = 8499.875
$ws.Range("K74").Value = 37506524
$ws.Range("L74").Value = 8499.875
$ws.Range("M74").Value = -37505588
$ws.Range("N74").Value = -10371.875

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 18757512
$ws.Range("I77").Value = 37506524
$ws.Range("J77").Value = 8499.875
$ws.Range("K77").Value = 187532620
$ws.Range("L77").Value = 42499.375
$ws.Range("M77").Value = -187527940
$ws.Range("N77").Value = -51859.375

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 18947.7
$ws.Range("I116").Value = 18854.572
$ws.Range("K116").Value = 18854.572
$ws.Range("M116").Value = -15412.572

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1515.2742
$ws.Range("I137").Value = 1467.5264
$ws.Range("J137").Value = 2059.6
$ws.Range("K137").Value = 4402.5792
$ws.Range("L137").Value = 6178.799999999999
$ws.Range("M137").Value = -1852.5792
$ws.Range("N137").Value = -11278.8

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6032.6665
$ws.Range("I132").Value = 3486.75
$ws.Range("J132").Value = 11124.5
$ws.Range("K132").Value = 10460.25
$ws.Range("L132").Value = 33373.5
$ws.Range("M132").Value = -7930.25
$ws.Range("N132").Value = -38433.5

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1799.4
$ws.Range("I99").Value = 1625
$ws.Range("J99").Value = 2497
$ws.Range("K99").Value = 1625
$ws.Range("L99").Value = 2497
$ws.Range("M99").Value = -127
$ws.Range("N99").Value = -5493

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4913.579
$ws.Range("I107").Value = 4913.579
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4913.579
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2993.579
$ws.Range("N107").ClearContents()

# CRP row 4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 6499.5
$ws.Range("J4").Value = 6499.5
$ws.Range("L4").Value = 6499.5
$ws.Range("N4").Value = -6723.5

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4007.6296
$ws.Range("I58").Value = 2850.5557
$ws.Range("J58").Value = 6321.778
$ws.Range("K58").Value = 2850.5557
$ws.Range("L58").Value = 6321.778
$ws.Range("M58").Value = -2647.5557
$ws.Range("N58").Value = -6727.778

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 5795.2144
$ws.Range("J94").Value = 5606.857
$ws.Range("L94").Value = 5606.857
$ws.Range("N94").Value = -6508.857

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2680
$ws.Range("I99").Value = 2666.6667
$ws.Range("J99").Value = 2800
$ws.Range("K99").Value = 2666.6667
$ws.Range("L99").Value = 2800
$ws.Range("M99").Value = -1168.6667
$ws.Range("N99").Value = -5796

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1243.2693
$ws.Range("I107").Value = 1014.86365
$ws.Range("J107").Value = 2499.5
$ws.Range("K107").Value = 1014.86365
$ws.Range("L107").Value = 2499.5
$ws.Range("M107").Value = 905.13635
$ws.Range("N107").Value = -6339.5

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2680
$ws.Range("I126").Value = 2666.6667
$ws.Range("J126").Value = 2800
$ws.Range("K126").Value = 8000.000100000001
$ws.Range("L126").Value = 8400
$ws.Range("M126").Value = -5530.000100000001
$ws.Range("N126").Value = -13340

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4007.6296
$ws.Range("I136").Value = 2850.5557
$ws.Range("J136").Value = 6321.778
$ws.Range("K136").Value = 8551.667099999999
$ws.Range("L136").Value = 18965.334
$ws.Range("M136").Value = -6001.667099999999
$ws.Range("N136").Value = -24065.334

# CUL row 8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 661.8333
$ws.Range("I8").Value = 661.8333
$ws.Range("K8").Value = 1985.4999
$ws.Range("M8").Value = -1846.4999

# CUL row 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 170.61111
$ws.Range("I11").Value = 173.25
$ws.Range("K11").Value = 519.75
$ws.Range("M11").Value = -379.75

# CUL row 76
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 13713.5
$ws.Range("J76").Value = 15773
$ws.Range("L76").Value = 47319
$ws.Range("N76").Value = -48085

# CUL row 79
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H79").Value = 13713.5
$ws.Range("J79").Value = 15773
$ws.Range("L79").Value = 47319
$ws.Range("N79").Value = -49971

# CUL row 106
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 16638
$ws.Range("J106").Value = 16638
$ws.Range("L106").Value = 49914
$ws.Range("N106").Value = -51806

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 905.2941
$ws.Range("I82").Value = 804.63635
$ws.Range("J82").Value = 1089.8334
$ws.Range("K82").Value = 804.63635
$ws.Range("L82").Value = 1089.8334
$ws.Range("M82").Value = -443.63635
$ws.Range("N82").Value = -1811.8334

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 905.2941
$ws.Range("I85").Value = 804.63635
$ws.Range("J85").Value = 1089.8334
$ws.Range("K85").Value = 804.63635
$ws.Range("L85").Value = 1089.8334
$ws.Range("M85").Value = 443.36365
$ws.Range("N85").Value = -3585.8334

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2538
$ws.Range("I93").Value = 2053.3076
$ws.Range("J93").Value = 3798.2
$ws.Range("K93").Value = 2053.3076
$ws.Range("L93").Value = 3798.2
$ws.Range("M93").Value = -805.3076000000001
$ws.Range("N93").Value = -6294.2

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5280.8887
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10916.667
$ws.Range("J62").Value = 12625
$ws.Range("L62").Value = 12625
$ws.Range("N62").Value = -13873

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 10916.667
$ws.Range("J65").Value = 12625
$ws.Range("L65").Value = 63125
$ws.Range("N65").Value = -69365

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1766.3846
$ws.Range("I81").Value = 1451.2727
$ws.Range("J81").Value = 3499.5
$ws.Range("K81").Value = 2902.5454
$ws.Range("L81").Value = 6999
$ws.Range("M81").Value = -1841.5454
$ws.Range("N81").Value = -9121

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1766.3846
$ws.Range("I84").Value = 1451.2727
$ws.Range("J84").Value = 3499.5
$ws.Range("K84").Value = 14512.727
$ws.Range("L84").Value = 34995
$ws.Range("M84").Value = -9208.726999999999
$ws.Range("N84").Value = -45603

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1668698.5
$ws.Range("I107").Value = 2859383.2
$ws.Range("J107").Value = 1739.8
$ws.Range("K107").Value = 8578149.600000001
$ws.Range("L107").Value = 5219.4
$ws.Range("M107").Value = -8576229.600000001
$ws.Range("N107").Value = -9059.4

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3463.9048
$ws.Range("I136").Value = 2670.6843
$ws.Range("J136").Value = 10999.5
$ws.Range("K136").Value = 8012.0529
$ws.Range("L136").Value = 32998.5
$ws.Range("M136").Value = -5462.0529
$ws.Range("N136").Value = -38098.5

Write-Output "Applied all Moogle_Profits updates"